$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Range("K3").Value = 2023
$ws.Range("K4").Value = 2245
$ws.Range("K5").Value = 344
$ws.Range("K6").Value = 1901
